$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new columns to make room for the new fields ---
# 1) Insert Round_4 and Round_5 before the old "total_score" column (G:H)
$ws.Range("G1:H1").EntireColumn.Insert()
# 2) Insert top_participant4 and top_participant5 before the old dummy columns (now N:O)
$ws.Range("N1:O1").EntireColumn.Insert()
# 3) Insert top_participant4_dummy and top_participant5_dummy before total_top_positions (now S:T)
$ws.Range("S1:T1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("G1").Value2 = "Round_4"
$ws.Range("H1").Value2 = "Round_5"
$ws.Range("I1").Value2 = "total_score"
$ws.Range("J1").Value2 = "average_score"
$ws.Range("K1").Value2 = "top_participant1"
$ws.Range("L1").Value2 = "top_participant2"
$ws.Range("M1").Value2 = "top_participant3"
$ws.Range("N1").Value2 = "top_participant4"
$ws.Range("O1").Value2 = "top_participant5"
$ws.Range("P1").Value2 = "top_participant1_dummy"
$ws.Range("Q1").Value2 = "top_participant2_dummy"
$ws.Range("R1").Value2 = "top_participant3_dummy"
$ws.Range("S1").Value2 = "top_participant4_dummy"
$ws.Range("T1").Value2 = "top_participant5_dummy"
$ws.Range("U1").Value2 = "total_top_positions"
$ws.Range("V1").Value2 = "final_evaluation"

# --- Numeric data columns D-J (rounds, total_score, average_score) ---
$ws.Range("D2").Value2 = 3
$ws.Range("D3").Value2 = 2
$ws.Range("D4").Value2 = 7
$ws.Range("D5").Value2 = 5
$ws.Range("D6").Value2 = 10
$ws.Range("E2").Value2 = 10
$ws.Range("E3").Value2 = 4
$ws.Range("E4").Value2 = 3
$ws.Range("E5").Value2 = 8
$ws.Range("E6").Value2 = 10
$ws.Range("F2").Value2 = 2
$ws.Range("F3").Value2 = 0
$ws.Range("F4").Value2 = 9
$ws.Range("F5").Value2 = 4
$ws.Range("F6").Value2 = 1
$ws.Range("G2").Value2 = 5
$ws.Range("G3").Value2 = 2
$ws.Range("G4").Value2 = 8
$ws.Range("G5").Value2 = 5
$ws.Range("G6").Value2 = 4
$ws.Range("H2").Value2 = 4
$ws.Range("H3").Value2 = 8
$ws.Range("H4").Value2 = 5
$ws.Range("H5").Value2 = 4
$ws.Range("H6").Value2 = 6
$ws.Range("I2").Value2 = 24
$ws.Range("I3").Value2 = 16
$ws.Range("I4").Value2 = 32
$ws.Range("I5").Value2 = 26
$ws.Range("I6").Value2 = 31
$ws.Range("J2").Value2 = 4.8
$ws.Range("J3").Value2 = 3.2
$ws.Range("J4").Value2 = 6.4
$ws.Range("J5").Value2 = 5.2
$ws.Range("J6").Value2 = 6.2

# --- Boolean columns K-O (top_participant1..5) ---
$ws.Range("K2").Value2 = $false
$ws.Range("K3").Value2 = $false
$ws.Range("K4").Value2 = $false
$ws.Range("K5").Value2 = $false
$ws.Range("K6").Value2 = $true
$ws.Range("L2").Value2 = $true
$ws.Range("L3").Value2 = $false
$ws.Range("L4").Value2 = $false
$ws.Range("L5").Value2 = $false
$ws.Range("L6").Value2 = $true
$ws.Range("M2").Value2 = $false
$ws.Range("M3").Value2 = $false
$ws.Range("M4").Value2 = $true
$ws.Range("M5").Value2 = $false
$ws.Range("M6").Value2 = $false
$ws.Range("N2").Value2 = $false
$ws.Range("N3").Value2 = $false
$ws.Range("N4").Value2 = $true
$ws.Range("N5").Value2 = $false
$ws.Range("N6").Value2 = $false
$ws.Range("O2").Value2 = $false
$ws.Range("O3").Value2 = $true
$ws.Range("O4").Value2 = $false
$ws.Range("O5").Value2 = $false
$ws.Range("O6").Value2 = $false

# --- Numeric columns P-V (dummy flags, total_top_positions, final_evaluation) ---
$ws.Range("P2").Value2 = 0
$ws.Range("P3").Value2 = 0
$ws.Range("P4").Value2 = 0
$ws.Range("P5").Value2 = 0
$ws.Range("P6").Value2 = 1
$ws.Range("Q2").Value2 = 1
$ws.Range("Q3").Value2 = 0
$ws.Range("Q4").Value2 = 0
$ws.Range("Q5").Value2 = 0
$ws.Range("Q6").Value2 = 1
$ws.Range("R2").Value2 = 0
$ws.Range("R3").Value2 = 0
$ws.Range("R4").Value2 = 1
$ws.Range("R5").Value2 = 0
$ws.Range("R6").Value2 = 0
$ws.Range("S2").Value2 = 0
$ws.Range("S3").Value2 = 0
$ws.Range("S4").Value2 = 1
$ws.Range("S5").Value2 = 0
$ws.Range("S6").Value2 = 0
$ws.Range("T2").Value2 = 0
$ws.Range("T3").Value2 = 1
$ws.Range("T4").Value2 = 0
$ws.Range("T5").Value2 = 0
$ws.Range("T6").Value2 = 0
$ws.Range("U2").Value2 = 1
$ws.Range("U3").Value2 = 1
$ws.Range("U4").Value2 = 2
$ws.Range("U5").Value2 = 0
$ws.Range("U6").Value2 = 2
$ws.Range("V2").Value2 = 0
$ws.Range("V3").Value2 = 0
$ws.Range("V4").Value2 = 1
$ws.Range("V5").Value2 = 1
$ws.Range("V6").Value2 = 1

Write-Host "done"
